# Uploader_WBS.xlsx: insert a new WBS row (row 5) for the GitHub repository
# setup task, pushing the existing rows (old 5-17) down to (6-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5 - shifts rows 5..17 down to 6..18
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with the new task info
$ws.Range("C5").Value = "GITHUB repository setup ."
$ws.Range("D5").Value = "InProgress"
$ws.Range("E5").Value = "25-07-2016"
